# Qualifier 2 DC vs KKR.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 ---

# Row 67 (Eliminator row) - Sushma's points for this match come in as 0.
$ws1.Range("Q67").Value = 0

# Row 68 (Qualifier 2: DC vs KKR) - match name + each player's raw points.
$ws1.Range("C68").Value = "DC vs KKR"
$ws1.Range("E68").Value = 80
$ws1.Range("H68").Value = 60
$ws1.Range("K68").Value = 20
$ws1.Range("N68").Value = 100
$ws1.Range("Q68").Value = 0
$ws1.Range("T68").Value = 40

# Rows 80-85 coin-split table, column E = Qualifier 2 coin contribution per player.
$ws1.Range("E80").Value = 5
$ws1.Range("E81").Value = 7
$ws1.Range("E82").Value = 0
$ws1.Range("E83").Value = 3
$ws1.Range("E85").Value = 11

# --- Sheet2 ---

# Qualifier 2 (DC vs KKR) rank predictions.
$ws2.Range("K44").Value = "Rapaka"
$ws2.Range("K45").Value = "Anantha"

# Sheet1 becomes the active/selected sheet in the saved workbook.
$ws1.Activate()
$ws1.Range("Q86").Select()

# Sheet2 keeps its own remembered selection for when it's revisited.
$ws2.Select()
$ws2.Range("H48").Select()

# Restore Sheet1 as the active sheet before saving.
$ws1.Activate()

$wb.Save()
